$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 85, shifting rows 85:89 down to 86:90
$ws.Rows.Item(85).Insert()

# Populate the new row 85 with data (copy constant columns from the row below, set new values)
$ws.Cells.Item(85, 1).Value = 2
$ws.Cells.Item(85, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(85, 3).Value = "Coquimbo"
$ws.Cells.Item(85, 4).Value = 45106
$ws.Cells.Item(85, 4).NumberFormat = $ws.Cells.Item(86, 4).NumberFormat
$ws.Cells.Item(85, 5).Value = 4
$ws.Cells.Item(85, 6).Value = 100112022
$ws.Cells.Item(85, 7).Value = "Arveja Verde"
$ws.Cells.Item(85, 8).Value = "Perfection"
$ws.Cells.Item(85, 9).Value = "Primera"
$ws.Cells.Item(85, 10).Value = 900
$ws.Cells.Item(85, 11).Value = 26000
$ws.Cells.Item(85, 12).Value = 28000
$ws.Cells.Item(85, 13).Value = 27000
$ws.Cells.Item(85, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(85, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(85, 16).Value = 1080
$ws.Cells.Item(85, 17).Value = 25
$ws.Cells.Item(85, 18).Value = "Hortaliza"
